$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added as row 103 in the source data.
# This pushes the existing rows 103-135 down to 104-136 (same data, same order).
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A103").Value = 10
$ws.Range("B103").Value = "Vega Modelo de Temuco"
$ws.Range("C103").Value = "La Araucanía"
$ws.Range("D103").Value = 45229
$ws.Range("E103").Value = 9
$ws.Range("F103").Value = 300000001
$ws.Range("G103").Value = "Rabanito"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 35
$ws.Range("K103").Value = 8000
$ws.Range("L103").Value = 8000
$ws.Range("M103").Value = 8000
$ws.Range("N103").Value = "$/docena de paquetes"
$ws.Range("O103").Value = "Provincia de Cautín"
$ws.Range("P103").Value = 667
$ws.Range("Q103").Value = 12
$ws.Range("R103").Value = "Hortaliza"
